# The deck ships two themes:
#   ppt/theme/theme1.xml -> "Office Theme" (used by the notes master)
#   ppt/theme/theme2.xml -> "Integral"     (used by the slide master)
# The commit swaps their contents: theme1.xml becomes the "Integral"
# theme and theme2.xml becomes the stock "Office Theme". The fonts and
# format (fill/line/effect) schemes are identical ("Office") between the
# two themes, so the only real difference is the 12-slot theme colour
# scheme. We reach the presentation's theme through the slide master and
# rewrite its colour scheme in place to the "Office Theme" palette so the
# live theme part ends up holding the colours the diff expects there.

function ToComRgb([int]$hexRRGGBB) {
    # PowerPoint's ColorFormat/ThemeColor .RGB uses the Win32 COLORREF
    # layout (0x00BBGGRR) rather than the usual 0xRRGGBB order, so the
    # bytes have to be swapped before assigning.
    $r = ($hexRRGGBB -shr 16) -band 0xFF
    $g = ($hexRRGGBB -shr 8) -band 0xFF
    $b = $hexRRGGBB -band 0xFF
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

# Slot order (1-based, matches ThemeColorScheme.Item index):
#  1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#  8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    0x000000, # dk1
    0xFFFFFF, # lt1
    0x44546A, # dk2
    0xE7E6E6, # lt2
    0x5B9BD5, # accent1
    0xED7D31, # accent2
    0xA5A5A5, # accent3
    0xFFC000, # accent4
    0x4472C4, # accent5
    0x70AD47, # accent6
    0x0563C1, # hlink
    0x954F72  # folHlink
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = ToComRgb $officeThemeColors[$i - 1]
}
